# Scen_RES_SHARE_90%_24_7.xlsx - "Add files via upload" edit
#
# Semantic content change: the LimType value in row 5 (column G) of the
# INS_1 sheet is changed from "FX" to "LO".
#
# The user also had cell J12 selected when the file was last saved (the
# previous selection was C5:D5), so we move the selection there too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS_1")

# Update the LimType cell (G5) from "FX" to "LO".
$ws.Range("G5").Value = "LO"

# Reflect the saved selection/active cell (J12) recorded in the sheet view.
$ws.Range("J12").Select()
